$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range('D2') '27.444.45'
Set-TextCell $ws.Range('E2') '  -2.04%  '
Set-TextCell $ws.Range('D3') '1.836.63'
Set-TextCell $ws.Range('E3') '  -2.56%  '
Set-TextCell $ws.Range('D4') '1.004'
Set-TextCell $ws.Range('E4') '  -0.71%  '
Set-TextCell $ws.Range('D5') '332.45'
Set-TextCell $ws.Range('E5') '  -1.02%  '
Set-TextCell $ws.Range('E6') '  -0.86%  '
Set-TextCell $ws.Range('D7') '0.4607'
Set-TextCell $ws.Range('E7') '  -3.26%  '
Set-TextCell $ws.Range('D8') '0.3819'
Set-TextCell $ws.Range('E8') '  -3.35%  '
Set-TextCell $ws.Range('D9') '46.38'
Set-TextCell $ws.Range('E9') '  -1.17%  '
Set-TextCell $ws.Range('D10') '0.07905'
Set-TextCell $ws.Range('E10') '  -1.58%  '
Set-TextCell $ws.Range('D11') '0.9746'
Set-TextCell $ws.Range('E11') '  -4.54%  '
Set-TextCell $ws.Range('D12') '21.14'
Set-TextCell $ws.Range('E12') '  -3.64%  '
Set-TextCell $ws.Range('B13') 'WrappedEther'
Set-TextCell $ws.Range('C13') 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell $ws.Range('D13') '1.871.01'
Set-TextCell $ws.Range('E13') '  -1.52%  '
Set-TextCell $ws.Range('B14') 'Polkadot'
Set-TextCell $ws.Range('C14') 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell $ws.Range('D14') '5.910'
Set-TextCell $ws.Range('E14') '  -2.57%  '
Set-TextCell $ws.Range('D15') '7.047'
Set-TextCell $ws.Range('E15') '  -2.25%  '
Set-TextCell $ws.Range('E16') '  -1.15%  '
Set-TextCell $ws.Range('D17') '87.96'
Set-TextCell $ws.Range('E17') '  -0.65%  '
Set-TextCell $ws.Range('D18') '0.06614'
Set-TextCell $ws.Range('E18') '  -1.78%  '
Set-TextCell $ws.Range('E19') '  -2.27%  '
Set-TextCell $ws.Range('D20') '17.05'
Set-TextCell $ws.Range('E20') '  -0.24%  '
Set-TextCell $ws.Range('E21') '  -0.84%  '
Set-TextCell $ws.Range('D22') '27.434.69'
Set-TextCell $ws.Range('D23') '5.365'
Set-TextCell $ws.Range('E23') '  -2.65%  '
Set-TextCell $ws.Range('D24') '10.84'
Set-TextCell $ws.Range('E24') '  -1.50%  '
Set-TextCell $ws.Range('D25') '2.301'
Set-TextCell $ws.Range('E25') '  -2.10%  '
Set-TextCell $ws.Range('D26') '157.14'
Set-TextCell $ws.Range('E26') '  -1.27%  '
Set-TextCell $ws.Range('D27') '19.41'
Set-TextCell $ws.Range('E27') '  -2.67%  '
Set-TextCell $ws.Range('D28') '2.070'
Set-TextCell $ws.Range('E28') '  -1.95%  '
Set-TextCell $ws.Range('D29') '5.333'
Set-TextCell $ws.Range('E29') '  -3.42%  '
Set-TextCell $ws.Range('D30') '118.95'
Set-TextCell $ws.Range('E30') '  -2.30%  '
Set-TextCell $ws.Range('D31') '0.9552'
Set-TextCell $ws.Range('E31') '  -2.65%  '
Set-TextCell $ws.Range('D32') '0.09293'
Set-TextCell $ws.Range('E32') '  -3.02%  '
Set-TextCell $ws.Range('D33') '3.564'
Set-TextCell $ws.Range('E33') '  -2.08%  '
Set-TextCell $ws.Range('E34') '  -1.96%  '
Set-TextCell $ws.Range('E35') '  -3.43%  '
Set-TextCell $ws.Range('D36') '0.05940'
Set-TextCell $ws.Range('E36') '  -2.38%  '
Set-TextCell $ws.Range('D37') '0.02193'
Set-TextCell $ws.Range('E37') '  -2.69%  '
Set-TextCell $ws.Range('D38') '8.072'
Set-TextCell $ws.Range('E38') '  -1.81%  '
Set-TextCell $ws.Range('E39') '  -4.08%  '
Set-TextCell $ws.Range('D40') '0.5802'
Set-TextCell $ws.Range('E40') '  -3.19%  '
Set-TextCell $ws.Range('D41') '0.1843'
Set-TextCell $ws.Range('E41') '  -2.74%  '
Set-TextCell $ws.Range('D42') '10.02'
Set-TextCell $ws.Range('E42') '  -3.31%  '
Set-TextCell $ws.Range('D43') '1.260'
Set-TextCell $ws.Range('E43') '  -1.55%  '
Set-TextCell $ws.Range('D44') '0.5493'
Set-TextCell $ws.Range('E44') '  -3.22%  '
Set-TextCell $ws.Range('D45') '11.96'
Set-TextCell $ws.Range('E45') '  -2.45%  '
Set-TextCell $ws.Range('D46') '1.870'
Set-TextCell $ws.Range('E46') '  -3.47%  '
Set-TextCell $ws.Range('D47') '0.06650'
Set-TextCell $ws.Range('E47') '  -2.25%  '
Set-TextCell $ws.Range('D48') '110.02'
Set-TextCell $ws.Range('E48') '  -2.39%  '
Set-TextCell $ws.Range('D49') '1.043'
Set-TextCell $ws.Range('E49') '  -2.83%  '
Set-TextCell $ws.Range('B50') 'PaxDollar'
Set-TextCell $ws.Range('C50') 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextCell $ws.Range('D50') '1.001'
Set-TextCell $ws.Range('E50') '  -0.98%  '
Set-TextCell $ws.Range('B51') 'BabyDogeCoin'
Set-TextCell $ws.Range('C51') 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextCell $ws.Range('D51') '0.00000000288'
Set-TextCell $ws.Range('E51') '  -1.49%  '
